# Updated symbol list on Mon Dec 19 19:42:08 UTC 2022 with GitHub Actions
# Refresh prices/volumes and re-rank a handful of coins on Sheet1.
# Price cells (column D) are text in the source data (e.g. "5.250"),
# so we lead with an apostrophe to force text entry and keep exact
# formatting (trailing zeros / no scientific notation) instead of
# letting Excel auto-convert the numeric-looking string to a Double.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'242.15"
$ws.Range("D4").Value = "'5.250"
$ws.Range("D5").Value = "'0.05600"
$ws.Range("D6").Value = "'3.375"
$ws.Range("D7").Value = "'6.378"
$ws.Range("D8").Value = "'0.8075"
$ws.Range("D9").Value = "'0.9167"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1426"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07293"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03219"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03025"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09276"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.611"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001651"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04699"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005807"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "'0.006366"
$ws.Range("D20").Value = "'0.004981"
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("D23").Value = "'0.0003104"
$ws.Range("D40").Value = "'0.03901"
$ws.Range("D41").Value = "'0.006973"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1033"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002913"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.007496"
$ws.Range("D45").Value = "'0.00005953"
$ws.Range("D47").Value = "'0.0005507"
$ws.Range("E47").Value = "46ACDXExchangeACXT"
$ws.Range("D48").Value = "'0.6833"
$ws.Range("D49").Value = "'0.06232"
$ws.Range("E49").Value = "48BOLOBOLOBestin24h"
$ws.Range("D51").Value = "'0.01011"
